# Updated symbol list on Fri Jan  6 19:35:33 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row : D value, E value ($null means "no change for that column")
$updates = @(
    @{Row=2;  D="258.93";        E="0.53%"}
    @{Row=3;  D="26.92";         E="-1.85%"}
    @{Row=4;  D="4.685";         E="2.49%"}
    @{Row=5;  D="0.06004";       E="2.05%"}
    @{Row=6;  D="6.667";         E="0.51%"}
    @{Row=7;  D="0.8580";        E="-0.08%"}
    @{Row=8;  D="0.9312";        E="0.57%"}
    @{Row=9;  D="0.1394";        E="-1.14%"}
    @{Row=10; D="0.04890";       E="35.61%"}
    @{Row=11; D=$null;           E="-0.95%"}
    @{Row=12; D="0.03132";       E="-2.94%"}
    @{Row=13; D="0.09139";       E="-0.45%"}
    @{Row=14; D="0.001529";      E="-1.38%"}
    @{Row=15; D="0.0006060";     E="0.00%"}
    @{Row=16; D="0.006161";      E="1.16%"}
    @{Row=17; D=$null;           E="-1.52%"}
    @{Row=19; D="2.166";         E="-1.74%"}
    @{Row=20; D="0.3111";        E="0.18%"}
    @{Row=21; D=$null;           E="1.52%"}
    @{Row=22; D="4.118";         E="6.61%"}
    @{Row=23; D="0.04224";       E="0.23%"}
    @{Row=24; D=$null;           E="-0.73%"}
    @{Row=25; D="0.004034";      E="-6.21%"}
    @{Row=26; D="0.0001200";     E="-0.01%"}
    @{Row=27; D=$null;           E="13.52%"}
    @{Row=40; D="0.03843";       E="0.23%"}
    @{Row=41; D="0.1114";        E="0.95%"}
    @{Row=42; D="0.003863";      E="-2.83%"}
    @{Row=43; D="0.002419";      E="0.41%"}
    @{Row=44; D="0.01525";       E="29.40%"}
    @{Row=45; D="0.00005114";    E="-6.40%"}
    @{Row=46; D="0.00000000750"; E="-0.04%"}
    @{Row=47; D=$null;           E="-16.74%"}
    @{Row=48; D="0.1303";        E="5.47%"}
    @{Row=49; D=$null;           E="-0.04%"}
    @{Row=50; D=$null;           E="-0.04%"}
)

foreach ($u in $updates) {
    $r = $u.Row
    if ($null -ne $u.D) {
        $ws.Cells.Item($r, 4).Value = "'" + $u.D
    }
    if ($null -ne $u.E) {
        $ws.Cells.Item($r, 5).Value = "'" + $u.E
    }
}
